# Refresh the crypto price/volume snapshot (scheduled GitHub Actions update).
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Column A (rank index) is
# untouched. Rows 44/45 additionally swap the Coin/Link pair (EnergySwap
# moves up to rank 42 / row 44, Frax drops to rank 43 / row 45).
#
# The D-column "price" values are free-form text in this sheet (several even
# contain multiple '.' separators, e.g. "24.881.40"), not real numbers. A
# plain `Range.Value = "49.60"` assignment would let Excel's type-inference
# re-parse it into the double 49.6 (dropping the trailing zero) or throw on
# the multi-dot ones. Forcing NumberFormat "@" (Text) before the write keeps
# the literal string, and resetting the Style back to "Normal" afterwards
# avoids leaving a stray text format behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($range, $text) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# row -> @{ B; C; D; E }  (omit a key to leave that column untouched)
$updates = [ordered]@{
    2  = @{ D = "24.881.40";   E = "  +2.13%  " }
    3  = @{ D = "1.711.02";    E = "  +1.90%  " }
    4  = @{                    E = "  -0.04%  " }
    5  = @{ D = "311.14";      E = "  +1.25%  " }
    6  = @{ D = "0.9975";      E = "  -0.04%  " }
    7  = @{ D = "0.3752";      E = "  +1.02%  " }
    8  = @{ D = "49.60";       E = "  +2.81%  " }
    9  = @{ D = "0.3450";      E = "  +0.26%  " }
    10 = @{ D = "1.209";       E = "  +2.10%  " }
    11 = @{ D = "0.07555";     E = "  +4.13%  " }
    12 = @{ D = "0.9988";      E = "  -0.12%  " }
    13 = @{ D = "21.14";       E = "  +3.62%  " }
    14 = @{ D = "6.323";       E = "  +3.53%  " }
    15 = @{ D = "7.051";       E = "  +4.38%  " }
    16 = @{ D = "1.711.28";    E = "  +2.03%  " }
    17 = @{ D = "0.00001137";  E = "  +2.50%  " }
    18 = @{ D = "0.06714";     E = "  -0.25%  " }
    19 = @{ D = "0.9974";      E = "  -0.10%  " }
    20 = @{ D = "85.03";       E = "  +4.71%  " }
    21 = @{ D = "17.35";       E = "  +5.42%  " }
    22 = @{ D = "6.395";       E = "  +4.90%  " }
    23 = @{ D = "13.18";       E = "  +10.32%  " }
    24 = @{ D = "24.870.97";   E = "  +2.29%  " }
    25 = @{ D = "2.452";       E = "  +0.99%  " }
    26 = @{ D = "2.803";       E = "  +5.30%  " }
    27 = @{ D = "20.45";       E = "  +4.48%  " }
    28 = @{ D = "151.74";      E = "  -0.29%  " }
    29 = @{ D = "132.34";      E = "  +3.99%  " }
    30 = @{ D = "1.900.87";    E = "  +2.10%  " }
    31 = @{ D = "1.243";       E = "  +28.47%  " }
    32 = @{ D = "7.006";       E = "  +11.00%  " }
    33 = @{ D = "4.257";       E = "  +5.60%  " }
    34 = @{ D = "1.858";       E = "  +6.55%  " }
    35 = @{ D = "13.92";       E = "  +13.28%  " }
    36 = @{ D = "0.08849";     E = "  +4.44%  " }
    37 = @{ D = "5.632";       E = "  +5.45%  " }
    38 = @{ D = "0.06689";     E = "  +4.12%  " }
    39 = @{ D = "9.198";       E = "  +2.06%  " }
    40 = @{ D = "0.02417";     E = "  +3.60%  " }
    41 = @{ D = "0.2251";      E = "  +6.61%  " }
    42 = @{ D = "1.280";       E = "  +1.33%  " }
    43 = @{ D = "0.6489";      E = "  +5.17%  " }
    44 = @{ B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "14.08";  E = "  +8.08%  " }
    45 = @{ B = "Frax";       C = "https://coinranking.com/coin/KfWtaeV1W+frax-frax";     D = "0.9971"; E = "  -0.04%  " }
    46 = @{ D = "0.6188";      E = "  +4.16%  " }
    47 = @{ D = "3.829";       E = "  +1.35%  " }
    48 = @{ D = "2.145";       E = "  +5.91%  " }
    49 = @{ D = "130.55";      E = "  +2.76%  " }
    50 = @{ D = "0.07314";     E = "  +1.47%  " }
    51 = @{ D = "80.01";       E = "  +5.30%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.Contains("B")) { $ws.Range("B$row").Value = $vals.B }
    if ($vals.Contains("C")) { $ws.Range("C$row").Value = $vals.C }
    if ($vals.Contains("D")) { Set-PriceText "D$row" $vals.D }
    if ($vals.Contains("E")) { $ws.Range("E$row").Value = $vals.E }
}
